$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) - new date serial values
$ws.Range("D2").Value = 44327
$ws.Range("D3").Value = 44330
$ws.Range("D4").Value = 44302
$ws.Range("D5").Value = 44309
$ws.Range("D6").Value = 44306
$ws.Range("D7").Value = 44323
$ws.Range("D8").Value = 44322
$ws.Range("D9").Value = 44313
$ws.Range("D10").Value = 44316

# Column M (Volumen)
$ws.Range("M3").Value = 60
$ws.Range("M4").Value = 80
$ws.Range("M6").Value = 80
$ws.Range("M8").Value = 60
$ws.Range("M9").Value = 120
$ws.Range("M10").Value = 120

# Row 5 / Row 8 Q, S, T swap (unidad de comercialización, precio $/kg, kg/unidad)
$ws.Range("Q5").Value = "$/caja 14 kilos granel"
$ws.Range("S5").Value = 821
$ws.Range("T5").Value = 14

$ws.Range("Q8").Value = "$/caja 10 kilos empedrada"
$ws.Range("S8").Value = 11500
$ws.Range("T8").Value = 1
